# Updated template and required fields check.
# Append new dictionary rows (57-63) to the flag_map Sheet1 table, reflecting
# newly added "file_location" record identifier plus several new
# "missing_*_units" / "missing_conc_bound_type" hard-stop flags.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New rows appended after the existing last row (56) -------------------
# Row 57: new "file_location" record-identifier field.
$ws.Cells.Item(57, 1).Value = "file_location"
$ws.Cells.Item(57, 2).Value = "Where the file was processed from a `"local`" file directory or from `"clowder`" with accompanying Clowder file ID."
$ws.Cells.Item(57, 3).Value = "Record Identifier"

# Rows 58-59: field names entered as a pair first, then both definitions
# filled in afterwards (matches original authoring/shared-string order).
$ws.Cells.Item(58, 1).Value = "missing_dose_volume_units"
$ws.Cells.Item(59, 1).Value = "missing_dermal_applied_area_units"
$ws.Cells.Item(58, 2).Value = "Missing dose_volume units with dose_volume entry"
$ws.Cells.Item(59, 2).Value = "Missing dermal_applied_area units with dermal_applied_area entry"
$ws.Cells.Item(58, 3).Value = "Hard Stop (Missing Required)"
$ws.Cells.Item(59, 3).Value = "Hard Stop (Missing Required)"

# Rows 60-63: remaining new hard-stop "missing *_units" flags, entered row by row.
$rows = @(
    @("missing_aerosol_particle_density_units", "Missing aerosol_particle_density units with aerosol_particle_density entry", "Hard Stop (Missing Required)"),
    @("missing_age_units", "Missing age units with age entry", "Hard Stop (Missing Required)"),
    @("missing_height_units", "Missing height units with height entry", "Hard Stop (Missing Required)"),
    @("missing_conc_bound_type", "Missing conc_bound_type when a bound entry is present", "Hard Stop (Missing Required)")
)

$startRow = 60
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$lastRow = $startRow + $rows.Count - 1

# NOTE: the existing AutoFilter / "_xlnm._FilterDatabase" defined name
# intentionally stay pinned to the original A1:C55 range (unchanged by this
# commit), so we do not touch AutoFilter() here.

# --- Update window/selection to mirror the author re-saving near the bottom
$ws.Range("C$lastRow").Select()
$excel.ActiveWindow.ScrollRow = $lastRow - 11

Write-Output "Appended rows through row $lastRow"
